# Work Profile and new tenant support
# Appends new certificate-registration-history rows to the "AMSIN" and
# "AMS" worksheets, and normalizes the formatting of AMS!A32:G32 to match
# the rest of the table.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Helper data: each row is Date(text), RunTime(serial), SprintName(text),
# TotalCases, PassCases, FailCases, TimeTaken
# ----------------------------------------------------------------------

# ---------------------------- AMSIN sheet ------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$amsinRows = @(
    @("79", "2023-03-09", 44994.56250695602, "cert174fstcycle", 51, 51, 0, 1.29),
    @("80", "2023-03-13", 44998.44736513889, "174certiflow",    51, 51, 0, 1.07),
    @("81", "2023-03-30", 45015.70042503472, "175sccert",       51, 50, 1, 1.15),
    @("82", "2023-03-31", 45016.50412170139, "175fnlcert",      51, 50, 1, 0.91),
    @("83", "2023-04-06", 45022.63169266203, "176newcert",      51, 48, 3, 3.34),
    @("84", "2023-04-07", 45023.66425609971, "176fstcert",      51, 50, 1, 1.09)
)

foreach ($row in $amsinRows) {
    $r = [int]$row[0]
    $dateText = $row[1]
    $runTime = $row[2]
    $sprintName = $row[3]
    $total = $row[4]
    $pass = $row[5]
    $fail = $row[6]
    $taken = $row[7]

    # The run-time column (B) always keeps the same date/time display
    # format as the row above it, so carry that formatting down first.
    $wsAmsin.Range("B" + ($r - 1)).Copy()
    $wsAmsin.Range("B" + $r).PasteSpecial(-4122)

    # Column A: store as literal text (leading apostrophe forces text so
    # the date-like string isn't auto-converted into a date serial).
    $wsAmsin.Cells.Item($r, 1).Value = "'" + $dateText

    # Column B: run time, numeric date/time serial value.
    $wsAmsin.Cells.Item($r, 2).Value = $runTime

    # Column C: sprint / build name (plain text).
    $wsAmsin.Cells.Item($r, 3).Value = $sprintName

    # Columns D-G: numeric counters / duration.
    $wsAmsin.Cells.Item($r, 4).Value = $total
    $wsAmsin.Cells.Item($r, 5).Value = $pass
    $wsAmsin.Cells.Item($r, 6).Value = $fail
    $wsAmsin.Cells.Item($r, 7).Value = $taken
}

# ------------------------------ AMS sheet -------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Normalize row 32's look to match the rest of the table (it was left
# with default/no explicit styling) and nudge its run-time value to the
# corrected figure.
$wsAms.Range("A31:G31").Copy()
$wsAms.Range("A32:G32").PasteSpecial(-4122)
$wsAms.Cells.Item(32, 2).Value = 44977.8228197338

$amsRows = @(
    @("33", "2023-03-07", 44992.75082172453, "173htfmarch",  51, 51, 0, 2.75),
    @("34", "2023-03-13", 44998.55603206019, "174betacert",  51, 51, 0, 1),
    @("35", "2023-03-13", 44998.85610415509, "174livecerrt", 51, 51, 0, 1.18),
    @("36", "2023-03-31", 45016.56446957176, "175btcert",    51, 51, 0, 0.97),
    @("37", "2023-03-31", 45016.8297484375,  "175certdev",   51, 48, 3, 1.52),
    @("38", "2023-03-31", 45016.84193472222, "175rrcer",     51, 50, 1, 1)
)

foreach ($row in $amsRows) {
    $r = [int]$row[0]
    $dateText = $row[1]
    $runTime = $row[2]
    $sprintName = $row[3]
    $total = $row[4]
    $pass = $row[5]
    $fail = $row[6]
    $taken = $row[7]

    $wsAms.Range("B" + ($r - 1)).Copy()
    $wsAms.Range("B" + $r).PasteSpecial(-4122)

    $wsAms.Cells.Item($r, 1).Value = "'" + $dateText
    $wsAms.Cells.Item($r, 2).Value = $runTime
    $wsAms.Cells.Item($r, 3).Value = $sprintName
    $wsAms.Cells.Item($r, 4).Value = $total
    $wsAms.Cells.Item($r, 5).Value = $pass
    $wsAms.Cells.Item($r, 6).Value = $fail
    $wsAms.Cells.Item($r, 7).Value = $taken
}
